# Adds a new "base_price" column (E) to the JPL sample data sheet,
# filling in a 25000 base price for every player row, matching the
# new "batch players upload" feature described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell E1: text + same formatting (bold/centered/wrap) as the
#     other header cells in row 1 ---
$ws.Range("E1").Value = "base_price"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data cells E2:E6: base price for each player ---
$ws.Range("E2").Value = 25000
$ws.Range("E3").Value = 25000
$ws.Range("E4").Value = 25000
$ws.Range("E5").Value = 25000
$ws.Range("E6").Value = 25000

# --- Column E width to fit the new data ---
$ws.Columns.Item(5).ColumnWidth = 9.92

# --- Header row needs to grow a bit to match the other rows ---
$ws.Rows.Item(1).RowHeight = 30

# --- Move/update the active selection to the last filled cell ---
$ws.Range("E6").Select() | Out-Null
